# Applies the LOQ4022.xlsx content update:
#  - Ativacao date 01/01/2012 -> 01/01/2022
#  - Updated "Objetivos" text (PT) and new "Objectives" text (EN)
#  - Updated "Programa resumido" text (PT) and new "Short syllabus" text (EN)
#  - Updated "Programa" text (PT) and new "Syllabus" text (EN)
#  - Updated "Metodo", "Criterio" and "Norma de recuperacao" texts
#  - Updated "Bibliografia" text
#
# NOTE: this runtime does not reliably bind user-defined function
# *parameters* (they come through as $null), so the helper routines below
# communicate via plain script-scope variables instead of parameters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteValues / xlPasteFormats paste-special constants
$xlPasteValues  = -4163
$xlPasteFormats = -4122

function Set-TextValue {
    # Uses $g_addr (cell address string) and $g_text (text to store).
    # Writing via a text formula + paste-values keeps the cell as a plain
    # shared string (t="s") even when the text looks like a date/number,
    # instead of Excel's usual auto-conversion to a numeric/date cell.
    $r = $ws.Range($g_addr)
    $r.Formula = '="' + $g_text + '"'
    $r.Copy() | Out-Null
    $r.PasteSpecial($xlPasteValues) | Out-Null
}

function Add-TranslationRow {
    # Uses $g_fmt1/$g_fmt2 (existing B/C cells to copy formatting from),
    # $g_tgt1/$g_tgt2 (new B/C cells to populate) and $g_text (shared text
    # for both target cells).
    $ws.Range($g_fmt1).Copy() | Out-Null
    $ws.Range($g_tgt1).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range($g_fmt2).Copy() | Out-Null
    $ws.Range($g_tgt2).PasteSpecial($xlPasteFormats) | Out-Null

    $g_addr = $g_tgt1
    Set-TextValue
    $g_addr = $g_tgt2
    Set-TextValue
}

# Ativacao: 01/01/2012 -> 01/01/2022 (both B8 and C8)
$g_addr = "B8"; $g_text = "01/01/2022"; Set-TextValue
$g_addr = "C8"; $g_text = "01/01/2022"; Set-TextValue

# Objetivos (PT) - updated wording
$g_addr = "B10"
$g_text = "Proporcionar aos alunos uma visão atual dos processos industriais que utilizam a conversão química como rota de transformação da matéria prima em produto. Serão estudados os processos das indústrias de química de base e de transformação."
Set-TextValue
$g_addr = "C10"
Set-TextValue

# Objectives (EN) - new row content (row 11 previously only had column A)
$g_fmt1 = "B10"; $g_fmt2 = "C10"; $g_tgt1 = "B11"; $g_tgt2 = "C11"
$g_text = "Provide students with a current view of industrial processes that use chemical conversion as a route to transform raw material into product. The processes of the basic chemical and transformation industries will be studied."
Add-TranslationRow

# Programa resumido (PT) - updated wording
$g_addr = "B14"
$g_text = "Introdução aos Processos Químicos Industriais; NPK / Fertilizantes; Ácido Sulfúrico; Cloro Álcalis; Papel e Celulose; Açúcar e álcool;  Processos Biotecnológicos;"
Set-TextValue
$g_addr = "C14"
Set-TextValue

# Short syllabus (EN) - new row content (row 15 previously only had column A)
$g_fmt1 = "B14"; $g_fmt2 = "C14"; $g_tgt1 = "B15"; $g_tgt2 = "C15"
$g_text = "Introduction to Industrial Chemical Processes; NPK / Fertilizers; Sulfuric Acid; Chlorine Alkali; Paper and Cellulose; Sugar and alcohol; Biotechnological Processes."
Add-TranslationRow

# Programa (PT) - updated wording
$g_addr = "B16"
$g_text = ".Introdução aos Processos Químicos Industriais; 2.NPK / Fertilizantes3.Ácido Sulfúrico; 4.Cloro Álcalis; 5.Papel e Celulose; 6.Açúcar e álcool; 7.Processos Biotecnológicos."
Set-TextValue
$g_addr = "C16"
Set-TextValue

# Syllabus (EN) - new row content (row 17 previously only had column A)
$g_fmt1 = "B16"; $g_fmt2 = "C16"; $g_tgt1 = "B17"; $g_tgt2 = "C17"
$g_text = "1. Introduction to Industrial Chemical Processes;2. NPK / Fertilizers3. Sulfuric Acid;4. Chlorine Alkali;5. Paper and Cellulose;6. Sugar and alcohol;7. Biotechnological Processes;"
Add-TranslationRow

# Metodo - updated wording
$g_addr = "B19"
$g_text = "Aulas expositivas, desenvolvimento de trabalhos e exercícios em sala e fora de sala de aula, discussão de casos práticos."
Set-TextValue
$g_addr = "C19"
Set-TextValue

# Criterio - updated wording
$g_addr = "B20"
$g_text = "Provas em sala, entrega de trabalhos e exercícios ou casos práticos elaborados fora de sala de aula."
Set-TextValue
$g_addr = "C20"
Set-TextValue

# Norma de recuperacao - updated wording
$g_addr = "B21"
$g_text = "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação."
Set-TextValue
$g_addr = "C21"
Set-TextValue

# Bibliografia - updated wording
$g_addr = "B22"
$g_text = "Ullmann’s encyclopedia of industrial chemistry; Editorial advisory board, Giuseppe Bellussi et al.; 7th, completely revised edition; Weinheim ; New York : WileyVCH, 2011.Encyclopedia of Chemical Processing; Edited by Sunggyu Lee; New York : Taylor & Francis, 2006.Kirk, Raymond Eller. Encyclopedia of chemical technology / Herman F.Mark et al. New York: John Wiley, 1984.Manual econômico da indústria química - MEIQ / Centro de Pesquisas e Desenvolvimento; 8ed; Camaçari: CEPED, 2007.Shreve, R. Norris; BRINK JR., J. A. Indústrias de processos químicos. Tradução de Horácio Macedo; 4.ed. Rio de Janeiro: Editora Guanabara Koogan, 2008, c1997.T.W. Graham Solomons, Craig B. Fryhle Hoboken, NJ. Organic chemistry; John Wiley, 9th ed; c2008.Revistas:Brazilian Journal of Chemical Engineering, São Paulo, SP: Brazilian Society of Chemical Engineering, v. 11, n. 1, 1995-;"
Set-TextValue
$g_addr = "C22"
Set-TextValue

Write-Output "edit complete"
